$wb = $excel.ActiveWorkbook

# --- About sheet: insert explanatory note about coal/gas after the
# "Some plant types are thus set to 1..." note (row 13), pushing the
# remaining notes down by two rows. ---
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Rows.Item(14).Insert()
$wsAbout.Rows.Item(14).Insert()
$wsAbout.Range("A14").Value = "However, we do not allow coal and gas to bid at their peak capacity factor to avoid overdispatch"
$wsAbout.Range("A15").Value = "of either type based on dispatch costs."

# --- BDSBaPCF sheet: coal (hard coal, lignite) and natural gas
# (nonpeaker) no longer bid at their peak capacity factor; they now
# bid at their expected capacity factor instead (flag set to 0). ---
$wsData = $wb.Worksheets.Item("BDSBaPCF")
$wsData.Range("B2").Value = 0   # hard coal
$wsData.Range("B3").Value = 0   # natural gas nonpeaker
$wsData.Range("B13").Value = 0  # lignite
